# Refined the code with possible error/logical error that may occur
# Add 5 new "Pet Category" asset rows to the Assets sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Assets")

$description = "Pet Category stores the Category of Pets"
$folder = "AssignmentFolder"
$names = @("Category Fish", "Category Cats", "Category Reptiles", "Category Birds", "Category Dogs")

$row = 3
foreach ($name in $names) {
    # Description first (column D), so the shared-string table picks up the
    # common description text before the per-row category name.
    $ws.Cells.Item($row, 4).Value = $description

    $ws.Cells.Item($row, 1).Value = $name
    $ws.Cells.Item($row, 1).Font.Name = "Segoe UI"

    $ws.Cells.Item($row, 2).Value = $name
    $ws.Cells.Item($row, 2).Font.Name = "Segoe UI"

    $ws.Cells.Item($row, 3).Value = $folder
    $ws.Cells.Item($row, 3).Font.Name = "Calibri"

    $ws.Cells.Item($row, 4).Font.Name = "Calibri"

    $row = $row + 1
}

# Reflect the author's final selection/scroll position on the Assets sheet.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("C11").Select()

Write-Output "Added pet category asset rows"
